$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended at the bottom of the user list: "hayet" / "988988".
# B8 must stay a text value (not be auto-coerced to a number), so force
# text formatting before assigning, then restore the default style so no
# stray formatting is left behind.
$ws.Range("A8").Value = "hayet"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "988988"
$ws.Range("B8").Style = "Normal"
